{"js": "// Merge the split runs of the document title (\"Flextable\" / \" \" / \"in\" / \" \" / \"word\")\n// into a single run reading \"Flextable in word\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The title is the first paragraph of the document (style \"Title\").\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.load(\"text,style\");\nawait context.sync();\n\n// Replacing the whole paragraph range's text collapses every run in the\n// paragraph into a single new run carrying the combined text.\nconst titleRange = titleParagraph.getRange(\"Whole\");\ntitleRange.insertText(\"Flextable in word\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Merge the split runs of the document title (\"Flextable\" / \" \" / \"in\" / \" \" / \"word\")\n# into a single run reading \"Flextable in word\".\n$d = $word.ActiveDocument\n\n# The title is the first paragraph of the document (style \"Title\").\n$titleParagraph = $d.Paragraphs(1)\n$titleRange = $titleParagraph.Range\n\n# A scoped find/replace across the whole paragraph range collapses every run\n# it spans into a single new run carrying the combined replacement text.\n$find = $titleRange.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Flextable in word\", $false, $false, $false, $false, $false, $true, 1, $false, \"Flextable in word\", 2)\n"}
